$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Remove the two extra "Straight Arrow Connector" shapes (ids 86, 87) ---
$s.Shapes.Item("Straight Arrow Connector 85").Delete()
$s.Shapes.Item("Straight Arrow Connector 86").Delete()

# --- Remove five of the translucent "Oval" decorations (ids 94, 95, 98, 99, 100) ---
$s.Shapes.Item("Oval 93").Delete()
$s.Shapes.Item("Oval 94").Delete()
$s.Shapes.Item("Oval 97").Delete()
$s.Shapes.Item("Oval 98").Delete()
$s.Shapes.Item("Oval 99").Delete()

# --- Remove the lone "A" textbox that used to sit with Group 105 (id 105) ---
$s.Shapes.Item("TextBox 104").Delete()

# --- Ungroup "Group 105" and keep only its first connector (id 107),
#     dropping its two sibling connectors (ids 108, 109) ---
$grp105 = $s.Shapes.Item("Group 105")
$ungrouped105 = $grp105.Ungroup()
$s.Shapes.Item("Straight Arrow Connector 107").Delete()
$s.Shapes.Item("Straight Arrow Connector 108").Delete()

# --- Remove four of the "B" textboxes (ids 112, 113, 114, 115) ---
$s.Shapes.Item("TextBox 111").Delete()
$s.Shapes.Item("TextBox 112").Delete()
$s.Shapes.Item("TextBox 113").Delete()
$s.Shapes.Item("TextBox 114").Delete()

# --- Ungroup "Group 115" and keep only its first connector (id 117),
#     dropping its two sibling connectors (ids 118, 119) ---
$grp115 = $s.Shapes.Item("Group 115")
$ungrouped115 = $grp115.Ungroup()
$s.Shapes.Item("Straight Arrow Connector 117").Delete()
$s.Shapes.Item("Straight Arrow Connector 118").Delete()

# --- Remove the "2.)" and "3.)" textboxes (ids 121, 122) ---
$s.Shapes.Item("TextBox 120").Delete()
$s.Shapes.Item("TextBox 121").Delete()
